$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.083.76"
$ws.Range("D3").Value = "2.977.52"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.560"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").Value = "3.451.10"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "2.973.61"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "52.149.13"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.177"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.74%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.84%  "
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").Value = "2.118.25"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.238"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0338"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.942"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
